# Commit: "Removed a few pictures"
#
# The first three slides of the deck (the ones built almost entirely
# around pictures/maps - "35m addresses", the lone picture slide, and
# the "82% of streets" slide) are removed. Deleting a slide in
# PowerPoint also removes its associated notes page, and the
# presentation's slide list / relationship bookkeeping is fixed up
# automatically.

$p = $ppt.ActivePresentation

# Delete the first three slides (always operate on index 1 - as each
# slide is removed, the next slide shifts up to index 1).
$p.Slides.Item(1).Delete()
$p.Slides.Item(1).Delete()
$p.Slides.Item(1).Delete()
